$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.661.04'
$ws.Range('E2').Value = '  +0.26%  '

$ws.Range('D3').Value = '1.698.20'
$ws.Range('E3').Value = '  +0.19%  '

$ws.Range('E4').Value = '  +0.15%  '

$origStyle_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.89'
$ws.Range('D5').Style = $origStyle_D5
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('E6').Value = '  +0.17%  '

$origStyle_D7 = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3924'
$ws.Range('D7').Style = $origStyle_D7
$ws.Range('E7').Value = '  -0.28%  '

$origStyle_D8 = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4041'
$ws.Range('D8').Style = $origStyle_D8
$ws.Range('E8').Value = '  +0.54%  '

$origStyle_D9 = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.519'
$ws.Range('D9').Style = $origStyle_D9
$ws.Range('E9').Value = '  -0.74%  '

$ws.Range('E10').Value = '  +0.18%  '

$origStyle_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '52.83'
$ws.Range('D11').Style = $origStyle_D11
$ws.Range('E11').Value = '  -1.86%  '

$origStyle_D12 = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08855'
$ws.Range('D12').Style = $origStyle_D12
$ws.Range('E12').Value = '  +1.06%  '

$origStyle_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.436'
$ws.Range('D13').Style = $origStyle_D13
$ws.Range('E13').Value = '  +3.05%  '

$ws.Range('E14').Value = '  +1.63%  '

$origStyle_D15 = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.152'
$ws.Range('D15').Style = $origStyle_D15
$ws.Range('E15').Value = '  +7.15%  '

$origStyle_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001320'
$ws.Range('D16').Style = $origStyle_D16
$ws.Range('E16').Value = '  -0.15%  '

$ws.Range('D17').Value = '1.701.55'
$ws.Range('E17').Value = '  +0.10%  '

$origStyle_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '99.52'
$ws.Range('D18').Style = $origStyle_D18
$ws.Range('E18').Value = '  -0.69%  '

$origStyle_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07026'
$ws.Range('D19').Style = $origStyle_D19
$ws.Range('E19').Value = '  -0.40%  '

$origStyle_D20 = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.74'
$ws.Range('D20').Style = $origStyle_D20
$ws.Range('E20').Value = '  +0.40%  '

$origStyle_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.076'
$ws.Range('D21').Style = $origStyle_D21
$ws.Range('E21').Value = '  +3.25%  '

$origStyle_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.006'
$ws.Range('D22').Style = $origStyle_D22
$ws.Range('E22').Value = '  +0.52%  '

$origStyle_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.62'
$ws.Range('D23').Style = $origStyle_D23
$ws.Range('E23').Value = '  +4.10%  '

$ws.Range('D24').Value = '24.652.14'
$ws.Range('E24').Value = '  +0.22%  '

$origStyle_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.136'
$ws.Range('D25').Style = $origStyle_D25
$ws.Range('E25').Value = '  +4.35%  '

$origStyle_D26 = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.354'
$ws.Range('D26').Style = $origStyle_D26
$ws.Range('E26').Value = '  +1.88%  '

$origStyle_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.65'
$ws.Range('D27').Style = $origStyle_D27
$ws.Range('E27').Value = '  +1.24%  '

$origStyle_D28 = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.02'
$ws.Range('D28').Style = $origStyle_D28
$ws.Range('E28').Value = '  +2.32%  '

$origStyle_D29 = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.684'
$ws.Range('D29').Style = $origStyle_D29
$ws.Range('E29').Value = '  +15.15%  '

$origStyle_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '135.63'
$ws.Range('D30').Style = $origStyle_D30
$ws.Range('E30').Value = '  +1.19%  '

$origStyle_D31 = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.159'
$ws.Range('D31').Style = $origStyle_D31
$ws.Range('E31').Value = '  -0.95%  '

$origStyle_D32 = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08960'
$ws.Range('D32').Style = $origStyle_D32
$ws.Range('E32').Value = '  +5.00%  '

$origStyle_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.600'
$ws.Range('D33').Style = $origStyle_D33
$ws.Range('E33').Value = '  +3.80%  '

$origStyle_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.068'
$ws.Range('D34').Style = $origStyle_D34
$ws.Range('E34').Value = '  -2.78%  '

$origStyle_D35 = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.973'
$ws.Range('D35').Style = $origStyle_D35
$ws.Range('E35').Value = '  +0.48%  '

$origStyle_D36 = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '11.08'
$ws.Range('D36').Style = $origStyle_D36
$ws.Range('E36').Value = '  -2.86%  '

$origStyle_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2753'
$ws.Range('D37').Style = $origStyle_D37
$ws.Range('E37').Value = '  +0.85%  '

$origStyle_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02850'
$ws.Range('D38').Style = $origStyle_D38
$ws.Range('E38').Value = '  +3.38%  '

$origStyle_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.43'
$ws.Range('D39').Style = $origStyle_D39
$ws.Range('E39').Value = '  -0.82%  '

$origStyle_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09136'
$ws.Range('D40').Style = $origStyle_D40
$ws.Range('E40').Value = '  +1.12%  '

$ws.Range('E41').Value = '  -0.64%  '

$origStyle_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7637'
$ws.Range('D42').Style = $origStyle_D42
$ws.Range('E42').Value = '  -0.61%  '

$origStyle_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '15.81'
$ws.Range('D43').Style = $origStyle_D43
$ws.Range('E43').Value = '  +3.18%  '

$origStyle_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7169'
$ws.Range('D44').Style = $origStyle_D44
$ws.Range('E44').Value = '  -0.18%  '

$origStyle_D45 = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.549'
$ws.Range('D45').Style = $origStyle_D45
$ws.Range('E45').Value = '  +1.40%  '

$origStyle_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.211'
$ws.Range('D46').Style = $origStyle_D46
$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('E47').Value = '  +0.23%  '

$ws.Range('E48').Value = '  -1.02%  '

$origStyle_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '139.94'
$ws.Range('D49').Style = $origStyle_D49
$ws.Range('E49').Value = '  -0.91%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$origStyle_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07968'
$ws.Range('D50').Style = $origStyle_D50
$ws.Range('E50').Value = '  -0.71%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$origStyle_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '90.30'
$ws.Range('D51').Style = $origStyle_D51
$ws.Range('E51').Value = '  +2.19%  '
